# "Added more user methods"
#
# On the "Methods" sheet, column A marks a row "Done" (green "Good" style)
# once that API method has been implemented. Ten more methods (rows 61-70,
# corresponding to the B61:B70 method names) have now been implemented, so
# mark each of their A-cells as "Done" using the same style already used
# for every other completed row (e.g. A55:A60).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Methods")
$ws.Activate()

for ($r = 61; $r -le 70; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = "Done"
    $cell.Style = "Good"
}

# Move the view/selection down to the newly completed rows, matching where
# the author was working (scrolled to A70).
$ws.Range("A70").Select()
